# Update the HK_R_acc_LG sheet: A1 keeps its label, A2:A49 get recomputed
# accuracy values (COBRAToolbox re-run produced new numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label in A1 is unchanged ("HK_R_acc_LG"); re-assert it for completeness.
$ws.Range("A1").Value = "HK_R_acc_LG"

$newValues = @(
    86.813186813186817,
    88.196988196988201,
    86.894586894586894,
    88.563288563288566,
    88.400488400488399,
    88.563288563288566,
    87.667887667887669,
    86.609686609686605,
    86.97598697598697,
    87.179487179487182,
    86.609686609686605,
    86.935286935286939,
    88.319088319088323,
    86.731786731786727,
    88.359788359788354,
    88.11558811558811,
    88.522588522588521,
    88.074888074888079,
    88.400488400488399,
    88.400488400488399,
    88.359788359788354,
    87.260887260887259,
    87.830687830687822,
    87.952787952787943,
    87.667887667887669,
    88.644688644688642,
    88.685388685388688,
    87.912087912087912,
    88.278388278388277,
    88.481888481888475,
    88.481888481888475,
    88.400488400488399,
    88.644688644688642,
    87.993487993487989,
    87.871387871387867,
    87.301587301587304,
    88.196988196988201,
    87.789987789987791,
    88.319088319088323,
    88.400488400488399,
    88.603988603988597,
    88.603988603988597,
    88.319088319088323,
    88.237688237688232,
    88.278388278388277,
    87.627187627187624,
    87.016687016687015,
    88.603988603988597
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
